$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 33 (shifts old rows 33.. down to 35..)
$ws.Rows.Item(33).Resize(2).Insert()

# New row 33 data
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = 44715
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = 100112036
$ws.Cells.Item(33, 7).Value = "Caigua"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 140
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 12).Value = 13000
$ws.Cells.Item(33, 13).Value = 12500
$ws.Cells.Item(33, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 625
$ws.Cells.Item(33, 17).Value = 20
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# New row 34 data
$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value = 44715
$ws.Cells.Item(34, 5).Value = 15
$ws.Cells.Item(34, 6).Value = 100112036
$ws.Cells.Item(34, 7).Value = "Caigua"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 150
$ws.Cells.Item(34, 11).Value = 9000
$ws.Cells.Item(34, 12).Value = 10000
$ws.Cells.Item(34, 13).Value = 9500
$ws.Cells.Item(34, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 475
$ws.Cells.Item(34, 17).Value = 20
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Match the date-format styling used in column D (numFmt "YYYY-MM-DD HH:MM:SS")
$ws.Range("D33:D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
